$d = $word.ActiveDocument

# --- 1) NAME line: "NAME: KAUSHIK NARAYANAN V" -> "NAME: " + "Ajay Kumar J" (two runs)
$rngName = $d.Content.Duplicate
$foundName = $rngName.Find.Execute("KAUSHIK NARAYANAN V", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundName) {
    $rngName.Text = "Ajay Kumar J"
    $rngName.Font.Bold = $true
    $rngName.Font.Name = "Times New Roman"
    $rngName.Font.NameBi = "Times New Roman"
}

# --- 2) REG NO line: "REG NO: 192321047" -> "REG NO: " + "192372052" (two runs)
$rngReg = $d.Content.Duplicate
$foundReg = $rngReg.Find.Execute("192321047", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundReg) {
    $rngReg.Text = "192372052"
    $rngReg.Font.Bold = $true
    $rngReg.Font.Name = "Times New Roman"
    $rngReg.Font.NameBi = "Times New Roman"
}

# --- 3) Merge the split "][j] == " / "nonTerminal" / ") {" runs (with proofErr wrapping
#        "nonTerminal") into a single run reading "][j] == nonTerminal) {"
$d.Content.Find.Execute("][j] == nonTerminal) {", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "][j] == nonTerminal) {", 2) | Out-Null

Write-Output "done"
